# Apply updated balance-sheet figures to the "IEX" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IEX")

# Row 4 - Inventory
$ws.Range("B4").Value = 305000000.0
$ws.Range("C4").Value = 290000000.0
$ws.Range("D4").Value = 302000000.0
$ws.Range("E4").Value = 325000000.0
$ws.Range("F4").Value = 341000000.0

# Row 12 - Accounts Payable
$ws.Range("B12").Value = 171000000.0
$ws.Range("C12").Value = 152000000.0
$ws.Range("D12").Value = 135000000.0
$ws.Range("E12").Value = 137000000.0
$ws.Range("F12").Value = 158000000.0

# Row 20 - Long Term Tax Liability (Deferred)
$ws.Range("B20").Value = 173000000.0
$ws.Range("C20").Value = 161000000.0
$ws.Range("D20").Value = 152000000.0
$ws.Range("E20").Value = 150000000.0
$ws.Range("F20").Value = 148000000.0

# Row 33 - Net Debt
$ws.Range("G33").Value = 216671000.0

# Row 34 - Total Debt
$ws.Range("G34").Value = 849252000.0
